$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Temporarily unprotect the sheet so values can be updated
$ws.Unprotect()

# Update the confidential disclaimer date in A18
$ws.Range("A18").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-03 for illustrative purposes only and are subject to change."

# Update D/E numeric values for rows 2-15
$ws.Range("D2").Value = 0.05772194610059739
$ws.Range("E2").Value = 0.001750345462920189

$ws.Range("D3").Value = 0.02354808135082754
$ws.Range("E3").Value = 0.01047223868800651

$ws.Range("D4").Value = 0.03169389490981903
$ws.Range("E4").Value = -0.0003775009437523291

$ws.Range("D5").Value = 0.03121930405864645
$ws.Range("E5").Value = 0.02548543689320404

$ws.Range("D6").Value = 0.03653950738187517
$ws.Range("E6").Value = 0.01506221349050429

$ws.Range("D7").Value = 0.01899545228966755
$ws.Range("E7").Value = 0.01211500683398081

$ws.Range("D8").Value = 0.004634240076155746
$ws.Range("E8").Value = -0.01419965576592075

$ws.Range("D9").Value = 0.006936736866718498
$ws.Range("E9").Value = 0.003832886163281035

$ws.Range("D10").Value = 0.07011182490431843
$ws.Range("E10").Value = 0.0130830489192264

$ws.Range("D11").Value = 0.07019158807258273
$ws.Range("E11").Value = 0.01249999999999996

$ws.Range("D12").Value = 0.1474448753088314
$ws.Range("E12").Value = -0.0008655510675128486

$ws.Range("D13").Value = 0.3864186508940142
$ws.Range("E13").Value = 0.001052723923151033

$ws.Range("D14").Value = 0.1145438977859459
$ws.Range("E14").Value = 0.0003655861564708651

$ws.Range("E15").Value = 0.003988304645690111

# Restore sheet protection as it was before the edit
$ws.Protect("D382")
